# Apply cryptos.xlsx price/volume updates (commit: Mon Oct  9 00:43:41 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.697.29"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "'1.620.27"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("D4").Value = "'0.995"
$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'209.91"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").Value = "'0.515"
$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").Value = "'23.15"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "'0.254"
$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("D11").Value = "'0.0876"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "'1.849.21"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "'1.599.57"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("D14").Value = "'3.99"
$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").Value = "'0.554"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").Value = "'64.78"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "'27.723.91"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "'226.34"
$ws.Range("E18").Value = "  -2.09%  "

$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").Value = "'0.0₃0714"
$ws.Range("E20").Value = "  -0.95%  "

$ws.Range("D21").Value = "'0.994"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").Value = "'4.30"
$ws.Range("E22").Value = "  -0.90%  "

$ws.Range("D23").Value = "'9.95"
$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("D24").Value = "'2.05"
$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("D25").Value = "'154.95"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("D28").Value = "'15.39"
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("D33").Value = "'3.08"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'1.397.53"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  +2.08%  "

$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  -1.06%  "

$ws.Range("D37").Value = "'2.33"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "'0.839"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("D41").Value = "'0.994"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").Value = "'1.84"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").Value = "'65.47"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").Value = "'5.39"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").Value = "'1.759.54"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("D47").Value = "'2.12"
$ws.Range("E47").Value = "  -3.01%  "

$ws.Range("D48").Value = "'87.67"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = "  +1.39%  "
